$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> becomes old row 3 values
$ws.Range("D2").Value = 44875
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 10

# Row 3 -> becomes old row 4 values
$ws.Range("D3").Value = 45222
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 1500

# Row 4 -> becomes old row 2 values
$ws.Range("D4").Value = 44855
$ws.Range("M4").Value = 25
$ws.Range("Q4").Value = '$/bandeja 5 kilos'
$ws.Range("R4").Value = 'Provincia de Los Andes'
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5
